$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expand used range dimension happens automatically when cells are written.

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Spp1"
$ws.Cells.Item(2,3).Value = "Itgb1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 12.486902
$ws.Cells.Item(2,8).Value = 37.460706
$ws.Cells.Item(2,9).Value = 0.01504353194025314
$ws.Cells.Item(2,10).Value = 0.01504353194025314
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 153.5290173333333
$ws.Cells.Item(2,14).Value = 460.587052
$ws.Cells.Item(2,15).Value = 0.3172206968818489
$ws.Cells.Item(2,16).Value = 0.317220696881849
$ws.Cells.Item(2,17).Value = 1917.101793597634
$ws.Cells.Item(2,18).Value = 17253.91614237871
$ws.Cells.Item(2,19).Value = 0.004772119685651453
$ws.Cells.Item(2,20).Value = 0.004772119685651454

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Spp1"
$ws.Cells.Item(3,3).Value = "Itgb1"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 12.486902
$ws.Cells.Item(3,8).Value = 37.460706
$ws.Cells.Item(3,9).Value = 0.01504353194025314
$ws.Cells.Item(3,10).Value = 0.01504353194025314
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 168.7997026666667
$ws.Cells.Item(3,14).Value = 506.3991080000001
$ws.Cells.Item(3,15).Value = 0.3487728915577651
$ws.Cells.Item(3,16).Value = 0.3487728915577651
$ws.Cells.Item(3,17).Value = 2107.785344827806
$ws.Cells.Item(3,18).Value = 18970.06810345025
$ws.Cells.Item(3,19).Value = 0.005246776134043683
$ws.Cells.Item(3,20).Value = 0.005246776134043684

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Spp1"
$ws.Cells.Item(4,3).Value = "Itgb1"
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 12.486902
$ws.Cells.Item(4,8).Value = 37.460706
$ws.Cells.Item(4,9).Value = 0.01504353194025314
$ws.Cells.Item(4,10).Value = 0.01504353194025314
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 68.09032333333333
$ws.Cells.Item(4,14).Value = 204.27097
$ws.Cells.Item(4,15).Value = 0.1406878008722904
$ws.Cells.Item(4,16).Value = 0.1406878008722904
$ws.Cells.Item(4,17).Value = 850.2371946116467
$ws.Cells.Item(4,18).Value = 7652.134751504821
$ws.Cells.Item(4,19).Value = 0.002116441426026273
$ws.Cells.Item(4,20).Value = 0.002116441426026274

# Row 5
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Spp1"
$ws.Cells.Item(5,3).Value = "Itgb1"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = 2
$ws.Cells.Item(5,6).Value = 0.6666666666666666
$ws.Cells.Item(5,7).Value = 12.486902
$ws.Cells.Item(5,8).Value = 37.460706
$ws.Cells.Item(5,9).Value = 0.01504353194025314
$ws.Cells.Item(5,10).Value = 0.01504353194025314
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 93.562673
$ws.Cells.Item(5,14).Value = 280.688019
$ws.Cells.Item(5,15).Value = 0.1933186106880956
$ws.Cells.Item(5,16).Value = 0.1933186106880956
$ws.Cells.Item(5,17).Value = 1168.307928609046
$ws.Cells.Item(5,18).Value = 10514.77135748142
$ws.Cells.Item(5,19).Value = 0.002908194694531728
$ws.Cells.Item(5,20).Value = 0.002908194694531729

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Spp1"
$ws.Cells.Item(6,3).Value = "Itgb1"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 46.08534733333334
$ws.Cells.Item(6,8).Value = 138.256042
$ws.Cells.Item(6,9).Value = 0.05552108878460485
$ws.Cells.Item(6,10).Value = 0.05552108878460485
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 153.5290173333333
$ws.Cells.Item(6,14).Value = 460.587052
$ws.Cells.Item(6,15).Value = 0.3172206968818489
$ws.Cells.Item(6,16).Value = 0.317220696881849
$ws.Cells.Item(6,17).Value = 7075.438089552021
$ws.Cells.Item(6,18).Value = 63678.94280596818
$ws.Cells.Item(6,19).Value = 0.01761243847589136
$ws.Cells.Item(6,20).Value = 0.01761243847589136

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Spp1"
$ws.Cells.Item(7,3).Value = "Itgb1"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 46.08534733333334
$ws.Cells.Item(7,8).Value = 138.256042
$ws.Cells.Item(7,9).Value = 0.05552108878460485
$ws.Cells.Item(7,10).Value = 0.05552108878460485
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 168.7997026666667
$ws.Cells.Item(7,14).Value = 506.3991080000001
$ws.Cells.Item(7,15).Value = 0.3487728915577651
$ws.Cells.Item(7,16).Value = 0.3487728915577651
$ws.Cells.Item(7,17).Value = 7779.192927156728
$ws.Cells.Item(7,18).Value = 70012.73634441054
$ws.Cells.Item(7,19).Value = 0.01936425067784204
$ws.Cells.Item(7,20).Value = 0.01936425067784204

# Row 8
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Spp1"
$ws.Cells.Item(8,3).Value = "Itgb1"
$ws.Cells.Item(8,4).Value = "M2"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 46.08534733333334
$ws.Cells.Item(8,8).Value = 138.256042
$ws.Cells.Item(8,9).Value = 0.05552108878460485
$ws.Cells.Item(8,10).Value = 0.05552108878460485
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 68.09032333333333
$ws.Cells.Item(8,14).Value = 204.27097
$ws.Cells.Item(8,15).Value = 0.1406878008722904
$ws.Cells.Item(8,16).Value = 0.1406878008722904
$ws.Cells.Item(8,17).Value = 3137.966200855638
$ws.Cells.Item(8,18).Value = 28241.69580770074
$ws.Cells.Item(8,19).Value = 0.007811139883141241
$ws.Cells.Item(8,20).Value = 0.007811139883141243

# Row 9
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Spp1"
$ws.Cells.Item(9,3).Value = "Itgb1"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 46.08534733333334
$ws.Cells.Item(9,8).Value = 138.256042
$ws.Cells.Item(9,9).Value = 0.05552108878460485
$ws.Cells.Item(9,10).Value = 0.05552108878460485
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 93.562673
$ws.Cells.Item(9,14).Value = 280.688019
$ws.Cells.Item(9,15).Value = 0.1933186106880956
$ws.Cells.Item(9,16).Value = 0.1933186106880956
$ws.Cells.Item(9,17).Value = 4311.868282640089
$ws.Cells.Item(9,18).Value = 38806.8145437608
$ws.Cells.Item(9,19).Value = 0.01073325974773022
$ws.Cells.Item(9,20).Value = 0.01073325974773022

# Row 10
$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(10,2).Value = "Spp1"
$ws.Cells.Item(10,3).Value = "Itgb1"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 666.4749603333333
$ws.Cells.Item(10,8).Value = 1999.424881
$ws.Cells.Item(10,9).Value = 0.8029323328679479
$ws.Cells.Item(10,10).Value = 0.8029323328679479
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 153.5290173333333
$ws.Cells.Item(10,14).Value = 460.587052
$ws.Cells.Item(10,15).Value = 0.3172206968818489
$ws.Cells.Item(10,16).Value = 0.317220696881849
$ws.Cells.Item(10,17).Value = 102323.245737249
$ws.Cells.Item(10,18).Value = 920909.2116352407
$ws.Cells.Item(10,19).Value = 0.2547067541813391
$ws.Cells.Item(10,20).Value = 0.2547067541813391

# Row 11
$ws.Cells.Item(11,1).Value = "M2"
$ws.Cells.Item(11,2).Value = "Spp1"
$ws.Cells.Item(11,3).Value = "Itgb1"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 666.4749603333333
$ws.Cells.Item(11,8).Value = 1999.424881
$ws.Cells.Item(11,9).Value = 0.8029323328679479
$ws.Cells.Item(11,10).Value = 0.8029323328679479
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 168.7997026666667
$ws.Cells.Item(11,14).Value = 506.3991080000001
$ws.Cells.Item(11,15).Value = 0.3487728915577651
$ws.Cells.Item(11,16).Value = 0.3487728915577651
$ws.Cells.Item(11,17).Value = 112500.7751390451
$ws.Cells.Item(11,18).Value = 1012506.976251406
$ws.Cells.Item(11,19).Value = 0.2800410314595761
$ws.Cells.Item(11,20).Value = 0.2800410314595761

# Row 12
$ws.Cells.Item(12,1).Value = "M2"
$ws.Cells.Item(12,2).Value = "Spp1"
$ws.Cells.Item(12,3).Value = "Itgb1"
$ws.Cells.Item(12,4).Value = "M2"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 666.4749603333333
$ws.Cells.Item(12,8).Value = 1999.424881
$ws.Cells.Item(12,9).Value = 0.8029323328679479
$ws.Cells.Item(12,10).Value = 0.8029323328679479
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 68.09032333333333
$ws.Cells.Item(12,14).Value = 204.27097
$ws.Cells.Item(12,15).Value = 0.1406878008722904
$ws.Cells.Item(12,16).Value = 0.1406878008722904
$ws.Cells.Item(12,17).Value = 45380.49554266717
$ws.Cells.Item(12,18).Value = 408424.4598840046
$ws.Cells.Item(12,19).Value = 0.1129627841604494
$ws.Cells.Item(12,20).Value = 0.1129627841604494

# Row 13
$ws.Cells.Item(13,1).Value = "M2"
$ws.Cells.Item(13,2).Value = "Spp1"
$ws.Cells.Item(13,3).Value = "Itgb1"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 666.4749603333333
$ws.Cells.Item(13,8).Value = 1999.424881
$ws.Cells.Item(13,9).Value = 0.8029323328679479
$ws.Cells.Item(13,10).Value = 0.8029323328679479
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 93.562673
$ws.Cells.Item(13,14).Value = 280.688019
$ws.Cells.Item(13,15).Value = 0.1933186106880956
$ws.Cells.Item(13,16).Value = 0.1933186106880956
$ws.Cells.Item(13,17).Value = 62357.17877635564
$ws.Cells.Item(13,18).Value = 561214.6089872007
$ws.Cells.Item(13,19).Value = 0.1552217630665832
$ws.Cells.Item(13,20).Value = 0.1552217630665832

# Row 14
$ws.Cells.Item(14,1).Value = "sCs"
$ws.Cells.Item(14,2).Value = "Spp1"
$ws.Cells.Item(14,3).Value = "Itgb1"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 105.0040076666667
$ws.Cells.Item(14,8).Value = 315.012023
$ws.Cells.Item(14,9).Value = 0.1265030464071941
$ws.Cells.Item(14,10).Value = 0.1265030464071941
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 153.5290173333333
$ws.Cells.Item(14,14).Value = 460.587052
$ws.Cells.Item(14,15).Value = 0.3172206968818489
$ws.Cells.Item(14,16).Value = 0.317220696881849
$ws.Cells.Item(14,17).Value = 16121.16211312513
$ws.Cells.Item(14,18).Value = 145090.4590181262
$ws.Cells.Item(14,19).Value = 0.040129384538967
$ws.Cells.Item(14,20).Value = 0.04012938453896701

# Row 15
$ws.Cells.Item(15,1).Value = "sCs"
$ws.Cells.Item(15,2).Value = "Spp1"
$ws.Cells.Item(15,3).Value = "Itgb1"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 105.0040076666667
$ws.Cells.Item(15,8).Value = 315.012023
$ws.Cells.Item(15,9).Value = 0.1265030464071941
$ws.Cells.Item(15,10).Value = 0.1265030464071941
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 168.7997026666667
$ws.Cells.Item(15,14).Value = 506.3991080000001
$ws.Cells.Item(15,15).Value = 0.3487728915577651
$ws.Cells.Item(15,16).Value = 0.3487728915577651
$ws.Cells.Item(15,17).Value = 17724.64527294172
$ws.Cells.Item(15,18).Value = 159521.8074564755
$ws.Cells.Item(15,19).Value = 0.04412083328630325
$ws.Cells.Item(15,20).Value = 0.04412083328630325

# Row 16
$ws.Cells.Item(16,1).Value = "sCs"
$ws.Cells.Item(16,2).Value = "Spp1"
$ws.Cells.Item(16,3).Value = "Itgb1"
$ws.Cells.Item(16,4).Value = "M2"
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 105.0040076666667
$ws.Cells.Item(16,8).Value = 315.012023
$ws.Cells.Item(16,9).Value = 0.1265030464071941
$ws.Cells.Item(16,10).Value = 0.1265030464071941
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 68.09032333333333
$ws.Cells.Item(16,14).Value = 204.27097
$ws.Cells.Item(16,15).Value = 0.1406878008722904
$ws.Cells.Item(16,16).Value = 0.1406878008722904
$ws.Cells.Item(16,17).Value = 7149.756833319145
$ws.Cells.Item(16,18).Value = 64347.81149987231
$ws.Cells.Item(16,19).Value = 0.01779743540267344
$ws.Cells.Item(16,20).Value = 0.01779743540267344

# Row 17
$ws.Cells.Item(17,1).Value = "sCs"
$ws.Cells.Item(17,2).Value = "Spp1"
$ws.Cells.Item(17,3).Value = "Itgb1"
$ws.Cells.Item(17,4).Value = "sCs"
$ws.Cells.Item(17,5).Value = 3
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 105.0040076666667
$ws.Cells.Item(17,8).Value = 315.012023
$ws.Cells.Item(17,9).Value = 0.1265030464071941
$ws.Cells.Item(17,10).Value = 0.1265030464071941
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 93.562673
$ws.Cells.Item(17,14).Value = 280.688019
$ws.Cells.Item(17,15).Value = 0.1933186106880956
$ws.Cells.Item(17,16).Value = 0.1933186106880956
$ws.Cells.Item(17,17).Value = 9824.455633005826
$ws.Cells.Item(17,18).Value = 88420.10069705243
$ws.Cells.Item(17,19).Value = 0.02445539317925046
$ws.Cells.Item(17,20).Value = 0.02445539317925046
